$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "The application requires JRE 7 or newer." -> "... JRE 8 or newer."
#    with the number split into its own run and a _GoBack bookmark placed
#    right after it (matches Word's "last edit position" bookmark behavior).
# ---------------------------------------------------------------------------
$jrePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "The application requires JRE*") {
        $jrePara = $cand
        break
    }
}

$jreFull = $jrePara.Range
$jreFind = $d.Range($jreFull.Start, $jreFull.End)
$jreFind.Find.Execute("7", $false, $false, $false, $false, $false, $true, 1, $false, "8", 2) | Out-Null

# Force a run boundary around the newly inserted "8" by briefly bookmarking
# it, then drop a real (hidden) _GoBack bookmark right after it.
$tmpRange = $d.Range($jreFind.Start, $jreFind.End)
$d.Bookmarks.Add("zzTempSplit", $tmpRange) | Out-Null
$d.Bookmarks.Item("zzTempSplit").Delete()
$goBackRange = $d.Range($jreFind.End, $jreFind.End)
$d.Bookmarks.Add("_GoBack", $goBackRange) | Out-Null

# ---------------------------------------------------------------------------
# 2) Merge the 4 runs of the "swing interface / console interface" bullet
#    into a single run (pure whitespace/run-boundary cleanup, text unchanged).
# ---------------------------------------------------------------------------
$swingPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "The application can be run either*") {
        $swingPara = $cand
        break
    }
}

$swingFull = $swingPara.Range
$swingText = $swingFull.Text.Substring(0, $swingFull.Text.Length - 1)
$swingContent = $d.Range($swingFull.Start, $swingFull.End - 1)
$swingContent.Delete()
$d.Range($swingFull.Start, $swingFull.Start).InsertAfter($swingText) | Out-Null

# ---------------------------------------------------------------------------
# 3) Merge the trailing 4 runs ("need " / "to supply " / "arguments..." /
#    ", the first argument...") of the Runner bullet into a single run,
#    leaving the leading bold ", you will " run untouched.
# ---------------------------------------------------------------------------
$runnerPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "*you will need*Run Configuration*") {
        $runnerPara = $cand
        break
    }
}

$runnerFull = $runnerPara.Range
$needFind = $d.Range($runnerFull.Start, $runnerFull.End)
$needFind.Find.Execute("need ", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null

$mergedTail = "need to supply arguments through the " + [char]0x201C + "Run Configuration" + [char]0x201D + " dialog, the first argument is the path to the input CSV file and the second argument is the path to the output CSV file."

$tailContent = $d.Range($needFind.Start, $runnerFull.End - 1)
$tailContent.Delete()
$d.Range($needFind.Start, $needFind.Start).InsertAfter($mergedTail) | Out-Null

Write-Output "Done"
